# non-removal-patterns.xlsx : "Tight" sheet — drop the "Rename method" rows
# (rows 20-23) entirely, leaving only the blank (style-only) E column cells
# that used to host the shared "=Cn/Dn" percentage formula, and move the
# active selection to A18 (the last cell the author touched before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tight")

# Fully clear A20:D23 (values + styles) so the cells disappear from the
# sheet entirely rather than lingering as empty-but-styled cells.
$ws.Range("A20:D23").Clear() | Out-Null

# The E column keeps its percentage style (s="3") but loses its formula —
# clear contents only so the number format survives.
$ws.Range("E20:E23").ClearContents() | Out-Null

# Reflect the author's final selection before saving.
$ws.Range("A18").Select() | Out-Null
